$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1): swap J1/K1 and introduce Auto_Assign_Delivery column ---
$ws.Range("J1").Value = "Payment_Method"
$ws.Range("K1").Value = "Auto_Assign_Delivery"

# --- Row 2: update existing sample row ---
$ws.Range("A2").Value = "117"
$ws.Range("E2").Value = "Khilgaon"
$ws.Range("J2").Value = "bKash"
$ws.Range("K2").Value = "Yes"

# --- Row 3: add a brand new sample row ---
$ws.Range("A3").Value = "117"
$ws.Range("B3").Value = "Sanjida"
$ws.Range("C3").Value = "01900000000"
$ws.Range("D3").Value = "Dhaka"
$ws.Range("E3").Value = "Adabor"
$ws.Range("F3").Value = "Ada"
$ws.Range("G3").Value = "Chocolates2"
$ws.Range("H3").Value = "2"
$ws.Range("I3").Value = "250"
$ws.Range("J3").Value = "Cash on Delivery"
$ws.Range("K3").Value = "No"
$ws.Range("L3").Value = "Bring Carefully2!"

# --- View state: scroll so column B is leftmost, select B4 ---
$ws.Range("B4").Select()
$excel.ActiveWindow.ScrollColumn = 2
